# Update view-count figures (column F) on the "展览" and "全部类型" sheets
# to reflect the regenerated data output at commit 456a3b4.

$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F8").Value = 4133
$ws1.Range("F10").Value = 4862
$ws1.Range("F11").Value = 540
$ws1.Range("F12").Value = 1219

$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F9").Value = 4133
$ws4.Range("F11").Value = 4862
$ws4.Range("F12").Value = 540
$ws4.Range("F13").Value = 1219
